$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the "Zustand" (status) of the FilmBewertung related rows to "Erledigt"
$ws.Range("B8").Value = "Erledigt"
$ws.Range("B13").Value = "Erledigt"
$ws.Range("B18").Value = "Erledigt"
$ws.Range("B22").Value = "Erledigt"
$ws.Range("B52").Value = "Erledigt"

# Assign "Charly" as the responsible person ("Bearbeiter") for the newly
# finished FilmBewertung work items
$ws.Range("C18").Value = "Charly"
$ws.Range("C52").Value = "Charly"

# Widen column A so the longer texts fit nicely
$ws.Columns.Item(1).ColumnWidth = 45.42578125

# Move/update the current selection to reflect where the user ended up editing
$ws.Range("C52").Select()
